{"js": "// The diff removes the very first paragraph of the document in its\n// entirety \u2014 a bold \"NAME: JEEVITHA.S-192324020\" line \u2014 leaving the\n// \"DATA BASE MANAGEMENT SYSTEM ...\" paragraph as the new first\n// paragraph. Everything else in the document is untouched.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph by its text (robust to any position/formatting\n// assumptions) and delete it outright (paragraph mark included).\nconst target = paragraphs.items.find((p) =>\n  p.text.includes(\"NAME: JEEVITHA.S-192324020\")\n);\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# The diff removes the very first paragraph of the document in its\n# entirety \u2014 a bold \"NAME: JEEVITHA.S-192324020\" line \u2014 leaving the\n# \"DATA BASE MANAGEMENT SYSTEM ...\" paragraph as the new first\n# paragraph. Everything else in the document is untouched.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*NAME: JEEVITHA.S-192324020*\") {\n        # Delete the whole paragraph, including its paragraph mark.\n        $p.Range.Delete()\n        break\n    }\n}\n"}
